$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 6 data rows (TW5/TW6/noSoln/noTW/1TW entries), leaving
# just the first three treatment rows.
$ws.Rows("5:10").Delete()

# Re-point the remaining rows at the correct treatment labels:
# row2 -> TW1-1, row3 stays TW2-1, row4 -> TW3-1
$ws.Range("B2").Value = "RPEmedia_TW1-1_d1_chamber1_channel1_5freq5sine_freq"
$ws.Range("B3").Value = "RPEmedia_TW2-1_d1_chamber1_channel1_5freq5sine_freq"
$ws.Range("B4").Value = "RPEmedia_TW3-1_d1_chamber1_channel1_5freq5sine_freq"

# Update the saved selection/active cell shown when the file is reopened.
$ws.Range("D12").Select() | Out-Null
